# Updates cryptos list values (prices, volume % changes) and corrects the
# Avalanche/Toncoin row ordering, matching the GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = '26.963.76'
    $ws.Range("D3").NumberFormat = "@"
    $ws.Range("D3").Value = '1.657.29'
    $ws.Range("E3").Value = '  +2.87%  '
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = '215.34'
    $ws.Range("E5").Value = '  +1.49%  '
    $ws.Range("E6").Value = '  +2.25%  '
    $ws.Range("E7").Value = '  -0.04%  '
    $ws.Range("D8").NumberFormat = "@"
    $ws.Range("D8").Value = '0.250'
    $ws.Range("E8").Value = '  +2.48%  '
    $ws.Range("E9").Value = '  +1.91%  '
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = '20.10'
    $ws.Range("E10").Value = '  +4.53%  '
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = '0.0881'
    $ws.Range("E11").Value = '  +4.20%  '
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = '1.889.69'
    $ws.Range("E12").Value = '  +2.80%  '
    $ws.Range("D13").NumberFormat = "@"
    $ws.Range("D13").Value = '1.658.20'
    $ws.Range("E13").Value = '  +2.80%  '
    $ws.Range("E14").Value = '  +2.02%  '
    $ws.Range("D15").NumberFormat = "@"
    $ws.Range("D15").Value = '0.522'
    $ws.Range("E15").Value = '  +2.57%  '
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = '65.50'
    $ws.Range("E16").Value = '  +3.02%  '
    $ws.Range("D17").NumberFormat = "@"
    $ws.Range("D17").Value = '26.973.75'
    $ws.Range("E17").Value = '  +2.12%  '
    $ws.Range("D18").NumberFormat = "@"
    $ws.Range("D18").Value = '236.04'
    $ws.Range("E18").Value = '  +0.21%  '
    $ws.Range("E19").Value = '  +1.66%  '
    $ws.Range("D20").NumberFormat = "@"
    $ws.Range("D20").Value = '7.75'
    $ws.Range("E20").Value = '  +1.15%  '
    $ws.Range("E21").Value = '  -0.04%  '
    $ws.Range("E22").Value = '  +3.96%  '
    $ws.Range("B23").Value = 'Toncoin'
    $ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = '2.23'
    $ws.Range("E23").Value = '  +2.19%  '
    $ws.Range("B24").Value = 'Avalanche'
    $ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    $ws.Range("D24").NumberFormat = "@"
    $ws.Range("D24").Value = '9.28'
    $ws.Range("E24").Value = '  +2.60%  '
    $ws.Range("D25").NumberFormat = "@"
    $ws.Range("D25").Value = '145.18'
    $ws.Range("E25").Value = '  -1.09%  '
    $ws.Range("E26").Value = '  +2.10%  '
    $ws.Range("E27").Value = '  +0.72%  '
    $ws.Range("D28").NumberFormat = "@"
    $ws.Range("D28").Value = '15.84'
    $ws.Range("E28").Value = '  +2.32%  '
    $ws.Range("E29").Value = '  -0.06%  '
    $ws.Range("E30").Value = '  +0.28%  '
    $ws.Range("E31").Value = '  +1.65%  '
    $ws.Range("D32").NumberFormat = "@"
    $ws.Range("D32").Value = '1.557.32'
    $ws.Range("E32").Value = '  +3.34%  '
    $ws.Range("E33").Value = '  +2.12%  '
    $ws.Range("D35").NumberFormat = "@"
    $ws.Range("D35").Value = '1.63'
    $ws.Range("E35").Value = '  +8.44%  '
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = '2.42'
    $ws.Range("E36").Value = '  -0.10%  '
    $ws.Range("D37").NumberFormat = "@"
    $ws.Range("D37").Value = '0.579'
    $ws.Range("E37").Value = '  +3.11%  '
    $ws.Range("D38").NumberFormat = "@"
    $ws.Range("D38").Value = '0.898'
    $ws.Range("E38").Value = '  +8.62%  '
    $ws.Range("E39").Value = '  +2.59%  '
    $ws.Range("E40").Value = '  +3.61%  '
    $ws.Range("E41").Value = '  -0.03%  '
    $ws.Range("D42").NumberFormat = "@"
    $ws.Range("D42").Value = '66.37'
    $ws.Range("E42").Value = '  +8.23%  '
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = '0.975'
    $ws.Range("E43").Value = '  +6.12%  '
    $ws.Range("E44").Value = '  +2.39%  '
    $ws.Range("D45").NumberFormat = "@"
    $ws.Range("D45").Value = '1.798.84'
    $ws.Range("E45").Value = '  +2.77%  '
    $ws.Range("E46").Value = '  +1.79%  '
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = '90.12'
    $ws.Range("E47").Value = '  +0.18%  '
    $ws.Range("E48").Value = '  +2.80%  '
    $ws.Range("D49").NumberFormat = "@"
    $ws.Range("D49").Value = '0.0998'
    $ws.Range("E49").Value = '  +4.16%  '
    $ws.Range("D50").NumberFormat = "@"
    $ws.Range("D50").Value = '0.0506'
    $ws.Range("E50").Value = '  +1.01%  '
    $ws.Range("D51").NumberFormat = "@"
    $ws.Range("D51").Value = '7.63'
    $ws.Range("E51").Value = '  +2.46%  '
